$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.851.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.344.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.666'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.91%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.66'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.34%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.87%  '

# Row 10
$ws.Range("E10").Value = '  -1.80%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.13%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.99%  '

# Row 13
$ws.Range("E13").Value = '  +0.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.59%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.694.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.52%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.30%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.901'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.40%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.343.52'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.72%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.801.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '78.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.70%  '

# Row 24
$ws.Range("E24").Value = '  +0.20%  '

# Row 25
$ws.Range("E25").Value = '  +2.53%  '

# Row 26
$ws.Range("E26").Value = '  +2.87%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.72%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.93%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.98%  '

# Row 32
$ws.Range("E32").Value = '  -1.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.133'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.47%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0746'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.34%  '

# Row 35
$ws.Range("E35").Value = '  -5.80%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.00%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.64%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.49%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +22.08%  '

# Row 41
$ws.Range("E41").Value = '  -3.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.18%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.45%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.106'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.42%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.195'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.76%  '

# Row 47
$ws.Range("E47").Value = '  +0.08%  '

# Row 48
$ws.Range("E48").Value = '  -3.22%  '

# Row 49
$ws.Range("B49").Value = 'BitTorrent-New'
$ws.Range("C49").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0155'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +53.90%  '

# Row 50
$ws.Range("E50").Value = '  -3.66%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.89%  '

